# "Fixed networking with less balls" - update mark-scheme input scores
# and let the dependent SUMPRODUCT/weighted-average formulas recalc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section "15" (rows 5-12): individual task marks out of 10 ---
$ws.Range("C10").Value = 9

# --- Section "9" (rows 15-18): individual task marks out of 10 ---
$ws.Range("C15").Value = 8
$ws.Range("C18").Value = 0

# --- Section "15" (rows 24-26) ---
$ws.Range("C25").Value = 10

# --- Section "15" (rows 37-41) ---
$ws.Range("C37").Value = 3
$ws.Range("C38").Value = 6
$ws.Range("C39").Value = 3
$ws.Range("C40").Value = 6
$ws.Range("C41").Value = 5

# --- Section "22" (rows 48-49) ---
$ws.Range("C48").Value = 3
$ws.Range("C49").Value = 3

# Recalculate the whole workbook so every SUMPRODUCT/weighted-average
# formula cell (C4, C14, C20, C23, C33, C36, C47, C51, C53:C56, ...)
# gets a refreshed cached value.
$excel.CalculateFullRebuild() | Out-Null

# --- Update the active/selected cell in the frozen bottom-right pane ---
$ws.Range("D39").Select() | Out-Null
